$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (2022-02-17 06:00)
$ws.Range("B7").Value = "[]"
$ws.Range("C7").Value = "[]"
$ws.Range("D7").Value = "[]"
$ws.Range("E7").Value = "[]"
$ws.Range("F7").Value = 0

# Row 8 (2022-02-17 07:00)
$ws.Range("B8").Value = "['Dacia SPRING']"
$ws.Range("C8").Value = "[0.3]"
$ws.Range("D8").Value = "[0.8500000000000002]"
$ws.Range("E8").Value = "[14.740000000000007]"
$ws.Range("F8").Value = 14.74000000000001

# Row 9 (2022-02-17 08:00)
$ws.Range("B9").Value = "['Opel CORSA', 'Others', 'SKODA ENYAQ 58kWh', 'VW ID.5']"
$ws.Range("C9").Value = "[0.3, 0.2, 0.2, 0.25]"
$ws.Range("D9").Value = "[0.8500000000000002, 0.8500000000000002, 0.7000000000000001, 0.7500000000000001]"
$ws.Range("E9").Value = "[24.75000000000001, 34.20986111111112, 29.0, 38.50000000000001]"
$ws.Range("F9").Value = 126.4598611111111

# Row 12 (2022-02-17 11:00)
$ws.Range("B12").Value = "['Tesla MODEL 3', 'Dacia SPRING']"
$ws.Range("C12").Value = "[0.35, 0.3]"
$ws.Range("D12").Value = "[0.9000000000000002, 0.8000000000000002]"
$ws.Range("E12").Value = "[27.500000000000014, 13.400000000000006]"
$ws.Range("F12").Value = 40.90000000000002

# Row 13 (2022-02-17 12:00)
$ws.Range("B13").Value = "['Audi Q4', 'Others', 'Fiat 500 E']"
$ws.Range("C13").Value = "[0.1, 0.35, 0.25]"
$ws.Range("D13").Value = "[0.8500000000000002, 0.9000000000000002, 0.9500000000000003]"
$ws.Range("E13").Value = "[57.45000000000001, 28.946805555555567, 16.800000000000008]"
$ws.Range("F13").Value = 103.1968055555556

# Row 14 (2022-02-17 13:00)
$ws.Range("B14").Value = "[]"
$ws.Range("C14").Value = "[]"
$ws.Range("D14").Value = "[]"
$ws.Range("E14").Value = "[]"
$ws.Range("F14").Value = 0

# Row 15 (2022-02-17 14:00)
$ws.Range("B15").Value = "['VW ID.3']"
$ws.Range("C15").Value = "[0.1]"
$ws.Range("D15").Value = "[0.8500000000000002]"
$ws.Range("E15").Value = "[43.500000000000014]"
$ws.Range("F15").Value = 43.50000000000001

# Row 16 (2022-02-17 15:00)
$ws.Range("B16").Value = "['Audi E-TRON', 'VW ID.4']"
$ws.Range("C16").Value = "[0.4, 0.1]"
$ws.Range("D16").Value = "[0.8000000000000002, 0.9500000000000003]"
$ws.Range("E16").Value = "[34.000000000000014, 65.45000000000002]"
$ws.Range("F16").Value = 99.45000000000003

# Row 17 (2022-02-17 16:00)
$ws.Range("B17").Value = "['Peugeot 208', 'Fiat 500 E']"
$ws.Range("C17").Value = "[0.2, 0.05]"
$ws.Range("D17").Value = "[0.8500000000000002, 0.7000000000000001]"
$ws.Range("E17").Value = "[29.250000000000007, 15.600000000000001]"
$ws.Range("F17").Value = 44.85000000000001

# Row 18 (2022-02-17 17:00)
$ws.Range("B18").Value = "[]"
$ws.Range("C18").Value = "[]"
$ws.Range("D18").Value = "[]"
$ws.Range("E18").Value = "[]"
$ws.Range("F18").Value = 0

# Row 32 (2022-02-18 07:00)
$ws.Range("B32").Value = "['Others']"
$ws.Range("C32").Value = "[0.35]"
$ws.Range("D32").Value = "[0.7500000000000001]"
$ws.Range("E32").Value = "[21.052222222222227]"
$ws.Range("F32").Value = 21.05222222222223

# Row 33 (2022-02-18 08:00)
$ws.Range("B33").Value = "['Tesla MODEL 3', 'Tesla MODEL 3', 'Tesla MODEL 3', 'MINI Cooper SE']"
$ws.Range("C33").Value = "[0.2, 0.35, 0.4, 0.2]"
$ws.Range("D33").Value = "[0.9500000000000003, 0.9500000000000003, 0.9000000000000002, 0.8500000000000002]"
$ws.Range("E33").Value = "[37.500000000000014, 30.000000000000014, 25.00000000000001, 18.785000000000004]"
$ws.Range("F33").Value = 111.2850000000001

# Row 36 (2022-02-18 11:00)
$ws.Range("B36").Value = "['Fiat 500 E', 'Others']"
$ws.Range("C36").Value = "[0.2, 0.25]"
$ws.Range("D36").Value = "[0.65, 0.9500000000000003]"
$ws.Range("E36").Value = "[10.8, 36.8413888888889]"
$ws.Range("F36").Value = 47.6413888888889

# Row 37 (2022-02-18 12:00)
$ws.Range("B37").Value = "['Smart FORTWO', 'Audi Q4']"
$ws.Range("C37").Value = "[0.4, 0.25]"
$ws.Range("D37").Value = "[0.9500000000000003, 0.9000000000000002]"
$ws.Range("E37").Value = "[9.680000000000005, 49.79000000000001]"
$ws.Range("F37").Value = 59.47000000000002

# Row 38 (2022-02-18 13:00)
$ws.Range("B38").Value = "['Others']"
$ws.Range("C38").Value = "[0.05]"
$ws.Range("D38").Value = "[0.65]"
$ws.Range("E38").Value = "[31.57833333333333]"
$ws.Range("F38").Value = 31.57833333333333

# Row 40 (2022-02-18 15:00)
$ws.Range("B40").Value = "['VW E-UP', 'Polestar 2', 'Others', 'Fiat 500 E']"
$ws.Range("C40").Value = "[0.3, 0.4, 0.2, 0.2]"
$ws.Range("D40").Value = "[0.9000000000000002, 0.8500000000000002, 0.8000000000000002, 0.8000000000000002]"
$ws.Range("E40").Value = "[22.08000000000001, 33.750000000000014, 31.578333333333337, 14.400000000000002]"
$ws.Range("F40").Value = 101.8083333333334

# Row 41 (2022-02-18 16:00)
$ws.Range("B41").Value = "['TESLA MODEL Y']"
$ws.Range("C41").Value = "[0.25]"
$ws.Range("D41").Value = "[0.65]"
$ws.Range("E41").Value = "[30.0]"
$ws.Range("F41").Value = 30
